$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.727.52"
$ws.Range("E2").Value = "  -5.68%  "

$ws.Range("D3").Value = "2.894.23"
$ws.Range("E3").Value = "  -3.97%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.43"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.20%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "2.889.75"
$ws.Range("E8").Value = "  -4.09%  "

$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.73"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000211"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.67%  "

$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "3.362.95"
$ws.Range("E16").Value = "  -4.10%  "

$ws.Range("D17").Value = "2.884.26"
$ws.Range("E17").Value = "  -4.34%  "

$ws.Range("D18").Value = "57.535.39"
$ws.Range("E18").Value = "  -5.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "407.70"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.655"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.61%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -2.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.98%  "

$ws.Range("E30").Value = "  -3.49%  "

$ws.Range("E31").Value = "  -3.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0953"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.904"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.97%  "

$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.27%  "

$ws.Range("D39").Value = "0.0₃0617"
$ws.Range("E39").Value = "  -9.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0343"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.59%  "

$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("D42").Value = "2.607.16"
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "356.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "

$ws.Range("E47").Value = "  -3.52%  "

$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.128"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.33%  "
